$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows and add new rows for the rescaled / time-incorporated data
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 7

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 4

# New rows 4 and 5: give A4/A5 the same style as A2/A3 (bold, bordered, centered)
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 2
